$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing rows (64-82) that held the longer sentences being dropped
$ws.Range("A64:A82").EntireRow.Delete()

# Rewrite each remaining row with the final (possibly re-ordered / re-worded) sentence
$ws.Range("A1").Value = '1|Imma give you guys a hint.|Imma give you guys a hint.'
$ws.Range("A2").Value = '2|This is a tennis ball cannon.|This is a tennis ball cannon.'
$ws.Range("A3").Value = '3|The only problem is, I''ve made some modifications to the cannon.|The only problem is, I''ve made some modifications to the cannon.'
$ws.Range("A4").Value = '4|Before we see if I actually live.|Before we see if I actually live.'
$ws.Range("A5").Value = '5|He sticks the landing like a pro just the same.|He sticks the landing like a pro just the same.'
$ws.Range("A6").Value = '6|Except if walk over here, we see number nine.|Except if walk over here, we see number nine.'
$ws.Range("A7").Value = '7|Chopsticks, play me a banger.|Chopsticks, play me a banger.'
$ws.Range("A8").Value = '8|This flashing light goes off, and everyone knows who not to high five.|This flashing light goes off, and everyone knows who not to high five.'
$ws.Range("A9").Value = '9|Honorable mention to the bathroom is the air freshener, which is actually filled with glitter bomb fart spray.|Honorable mention to the bathroom is the air freshener, which is actually filled with glitter bomb fart spray.'
$ws.Range("A10").Value = '10|Because subversion of expectations equals comedy.|Because subversion of expectations equals comedy.'
$ws.Range("A11").Value = '11|Coming in at number five is how we actully keep the lab location super secret.|Coming in at number five is how we actully keep the lab location super secret.'
$ws.Range("A12").Value = '12|If you were to walk in off the street or deliver a package here, this is what you''d see.|If you were to walk in off the street or deliver a package here, this is what you''d see.'
$ws.Range("A13").Value = '13|Clearly, whoever works here is not who you would expect would be hiding a high tech, state of the art Willy Wonka engineering lab.|Clearly, whoever works here is not who you would expect would be hiding a high tech, state of the art Willy Wonka engineering lab.'
$ws.Range("A14").Value = '17|And I brought my niece and nephew out here to tour the place.|And I brought my niece and nephew out here to tour the place.'
$ws.Range("A15").Value = '18|This was perhaps their favorite part.|This was perhaps their favorite part.'
$ws.Range("A16").Value = '19|Which just might have cemented my status in my ongoing quest to be the favorite uncle.|Which just might have cemented my status in my ongoing quest to be the favorite uncle.'
$ws.Range("A17").Value = '21|But after this, I''m gonna give this place a complete makeover.|But after this, I''m gonna give this place a complete makeover.'
$ws.Range("A18").Value = '22|So next time you find yourself waiting your turn at the barber.|So next time you find yourself waiting your turn at the barber.'
$ws.Range("A19").Value = '23|Or perhaps stuck in the waiting room at the dentist''s office.|Or perhaps stuck in the waiting room at the dentist''s office.'
$ws.Range("A20").Value = '24|You might want to get up and just start moving everything around on display.|You might want to get up and just start moving everything around on display.'
$ws.Range("A21").Value = '25|Because the entrance to the Crunch Lab just might be behind that oversized tooth decay poster hanging on the wall.|Because the entrance to the Crunch Lab just might be behind that oversized tooth decay poster hanging on the wall.'
$ws.Range("A22").Value = '27|And so for over a year, in addition to building this place, I''ve been working on a way to combine those three things.|And so for over a year, in addition to building this place, I''ve been working on a way to combine those three things.'
$ws.Range("A23").Value = '30|How do I foster that passion?|How do I foster that passion?'
$ws.Range("A24").Value = '31|Like what''s the first step?|Like what''s the first step?'
$ws.Range("A25").Value = '32|And I''ve just never had a great answer for that specific situation.|And I''ve just never had a great answer for that specific situation.'
$ws.Range("A26").Value = '33|Oh, and by the way, this is precisely where the toys are designed: right here in the Crunch Lab.|Oh, and by the way, this is precisely where the toys are designed: right here in the Crunch Lab.'
$ws.Range("A27").Value = '34|As you know, I love to take something humans do, and then build a robot to help me do it way better.|As you know, I love to take something humans do, and then build a robot to help me do it way better.'
$ws.Range("A28").Value = '35|Like with kicking field goals, skipping rocks, setting up dominoes, bowling, golfing, or even throwing bull''s eyes.|Like with kicking field goals, skipping rocks, setting up dominoes, bowling, golfing, or even throwing bull''s eyes.'
$ws.Range("A29").Value = '37|What kind of superpowers? Well I''m glad you asked.|What kind of superpowers? Well I''m glad you asked.'
$ws.Range("A30").Value = '38|We engineered the crud out of this thing to optimize six rapid-fire shots for maximum glide, with maximum accuracy.|We engineered the crud out of this thing to optimize six rapid-fire shots for maximum glide, with maximum accuracy.'
$ws.Range("A31").Value = '39|You will dominate any mere human who tries to challenge your creation.|You will dominate any mere human who tries to challenge your creation.'
$ws.Range("A32").Value = '42|It basically bumps into all these stationary air molecules and it causes it to deflect out.|It basically bumps into all these stationary air molecules and it causes it to deflect out.'
$ws.Range("A33").Value = '43|So for example, in the disk launcher, we learn all about how useful flywheels are.|So for example, in the disk launcher, we learn all about how useful flywheels are.'
$ws.Range("A34").Value = '46|And then in addition to the video, as backup, we''ve also got some simple booklet instructions.|And then in addition to the video, as backup, we''ve also got some simple booklet instructions.'
$ws.Range("A35").Value = '47|We maximized for sustainability in our material and packaging choices.|We maximized for sustainability in our material and packaging choices.'
$ws.Range("A36").Value = '48|Even the bags inside are compostable.|Even the bags inside are compostable.'
$ws.Range("A37").Value = '49|And the disk launcher is just the tip of the iceberg.|And the disk launcher is just the tip of the iceberg.'
$ws.Range("A38").Value = '50|In another box, we built this super cool Rube Goldberg catapult, and the gear badge here is for projectile motion.|In another box, we built this super cool Rube Goldberg catapult, and the gear badge here is for projectile motion.'
$ws.Range("A39").Value = '52|And there we learn about ratchets.|And there we learn about ratchets.'
$ws.Range("A40").Value = '53|And you get two of them, so you can compete in a Beyblades-style fight to the last coin standing.|And you get two of them, so you can compete in a Beyblades-style fight to the last coin standing.'
$ws.Range("A41").Value = '54|There''s a bunch more, but I want it to be a surprise for you, so rest assured: we designed our hearts out on every single last one.|There''s a bunch more, but I want it to be a surprise for you, so rest assured: we designed our hearts out on every single last one.'
$ws.Range("A42").Value = '55|And just like it says on the box, my goal with Crunch Labs is to help you think like an engineer.|And just like it says on the box, my goal with Crunch Labs is to help you think like an engineer.'
$ws.Range("A43").Value = '56|That means you think critically; you can break a problem down into manageable steps.|That means you think critically; you can break a problem down into manageable steps.'
$ws.Range("A44").Value = '58|We can crunch, and break, and fail things so we can learn fast.|We can crunch, and break, and fail things so we can learn fast.'
$ws.Range("A45").Value = '59|That''s literally the process of how I make every single one of my robots that eventually works out awesomely.|That''s literally the process of how I make every single one of my robots that eventually works out awesomely.'
$ws.Range("A46").Value = '60|There''s just no better teacher than a good crunch along the path to finding the best design.|There''s just no better teacher than a good crunch along the path to finding the best design.'
$ws.Range("A47").Value = '61|And any good engineer knows that''s true.|And any good engineer knows that''s true.'
$ws.Range("A48").Value = '63|Supplies are definitely limited on this first run, so if you don''t want to miss out, be sure to go check out our lovely website after this.|Supplies are definitely limited on this first run, so if you don''t want to miss out, be sure to go check out our lovely website after this.'
$ws.Range("A49").Value = '64|Speaking of which, I just gotta shout out Shopify, beause that''s how we built it.|Speaking of which, I just gotta shout out Shopify, beause that''s how we built it.'
$ws.Range("A50").Value = '65|They were super helpful, and I''ve been friends with them for a long time, and they give you all the tools to turn your idea into an actual business.|They were super helpful, and I''ve been friends with them for a long time, and they give you all the tools to turn your idea into an actual business.'
$ws.Range("A51").Value = '66|Alright, so continuing on with our top ten list of the coolest things here at Crunch Labs,|Alright, so continuing on with our top ten list of the coolest things here at Crunch Labs,'
$ws.Range("A52").Value = '67|At number four, if you look right here, we''ve got the start and finish line of the world''s longest Hot Wheels track at over a half a mile.|At number four, if you look right here, we''ve got the start and finish line of the world''s longest Hot Wheels track at over a half a mile.'
$ws.Range("A53").Value = '68|And set an official new world record.|And set an official new world record.'
$ws.Range("A54").Value = '71|Like how both ant colonies instinctively have made a graveyard in the corner for all their fallen comrades.|Like how both ant colonies instinctively have made a graveyard in the corner for all their fallen comrades.'
$ws.Range("A55").Value = '72|And now at number two, we''re back with the tennis ball cannon.|And now at number two, we''re back with the tennis ball cannon.'
$ws.Range("A56").Value = '73|But the cool thing here is, you might have noticed this glass suitcase attached to a chain, with ten thousand dollars cash in here.|But the cool thing here is, you might have noticed this glass suitcase attached to a chain, with ten thousand dollars cash in here.'
$ws.Range("A57").Value = '74|The deal is you can practice as much as you want, and I even set up some additional targets for that purpose.|The deal is you can practice as much as you want, and I even set up some additional targets for that purpose.'
$ws.Range("A58").Value = '75|But everyone gets only one official shot in their whole life to hit the bull''s eye and make it down the tube and into the box to unlock the ten thousand dollar suitcase.|But everyone gets only one official shot in their whole life to hit the bull''s eye and make it down the tube and into the box to unlock the ten thousand dollar suitcase.'
$ws.Range("A59").Value = '76|So for the avoidance of doubt, when you''re ready for your official shot, you just sign the wall and now the pressure''s really on.|So for the avoidance of doubt, when you''re ready for your official shot, you just sign the wall and now the pressure''s really on.'
$ws.Range("A60").Value = '77|And finally, the number one coolest thing about this place is you can actually come here.|And finally, the number one coolest thing about this place is you can actually come here.'
$ws.Range("A61").Value = '78|Every month, before we ship out the Crunch Lab build boxes, I essentially cover my eyes, throw a dart, and then I open whatever box I hit and set this inside.|Every month, before we ship out the Crunch Lab build boxes, I essentially cover my eyes, throw a dart, and then I open whatever box I hit and set this inside.'
$ws.Range("A62").Value = '80|And because I repeat this every month, that means you have a new chance to win every month.|And because I repeat this every month, that means you have a new chance to win every month.'
$ws.Range("A63").Value = '82|You can see a link to all the official rules in the video description.|You can see a link to all the official rules in the video description.'

# Restore the cursor/selection position recorded in the saved workbook
$null = $ws.Range("H15").Select()
